$d = $word.ActiveDocument

# ----------------------------------------------------------------------------
# Helper: replace the contents of a whole paragraph (but *not* its trailing
# paragraph mark) with a freshly built run, expressed as literal WordprocessingML.
#
# Using Range.InsertXML (rather than Range.Text = "..." or Find/Replace) gives
# full control over the resulting run: it lets a multi-run paragraph collapse
# into exactly one <w:r>, keeps the <w:t xml:space="preserve"> attribute that
# this document's runs otherwise always carry, and is immune to the
# smart-quote autocorrect that Find/Replace applies to straight quote
# characters.
# ----------------------------------------------------------------------------
function Set-ParagraphRun($paraIndex, $text) {
    $paraRange = $d.Paragraphs($paraIndex).Range
    # Exclude the final character (the paragraph mark) so InsertXML replaces
    # only the paragraph's visible content and keeps the paragraph itself.
    $target = $d.Range($paraRange.Start, $paraRange.End - 1)

    $escaped = $text -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'

    $runXml = '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:eastAsia="Arial"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="24"/><w:shd w:fill="auto" w:val="clear"/></w:rPr><w:t xml:space="preserve">' + $escaped + '</w:t></w:r>'

    $packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $runXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $target.InsertXML($packageXml)
}

# --- Paragraph 2 -------------------------------------------------------------
# Collapses the original 5 runs (plain / bold-space / plain / underlined "."
# / plain) describing Reinaldo Rogério Moreira and Print Express into a
# single new sentence about the company's "carro forte" service.
Set-ParagraphRun 2 'O "carro forte" da empresa é o serviço de encadernação em capa dura de Trabalhos de Conclusão de Curso dos alunos de Mestrado e Doutorado.'

# --- Paragraph 3 -------------------------------------------------------------
# "Nosso cliente ... Reinaldo ... Clayton, Douglas e Welington ... Douglas."
# becomes a more anonymized/generic description of the client and staff.
Set-ParagraphRun 3 'O cliente visa sucesso através da demanda que é feita pelos seus clientes, sendo xerox ou até mesmo as encadernações e também pela qualidade de seus produtos. O cliente é um dos operadores junto com três outros colaboradores que operam as máquinas de fotocópia em sua loja, e outro que produz as capas duras.'

# --- Paragraph 4 -------------------------------------------------------------
# "Uma das poucas interferências" -> "Uma das interferências".
Set-ParagraphRun 4 'Uma das interferências que ele tem com sua empresa no momento é a falta de marketing para a expansão do seu negócio. Tendo apenas uma faculdade em que atua, acaba não tendo um resultado tão bom quanto poderia ter, mas pensa que divulgar seu produto na internet seja uma boa estratégia, pois lá está seu público-alvo. '

# --- Paragraph 5 -------------------------------------------------------------
# "O cliente nos pede um sistema" -> "Ele pede um sistema".
Set-ParagraphRun 5 'Ele pede um sistema para conseguir gerenciar suas vendas e compras como uma planilha de gerenciamento. Um website onde as pessoas possam ver e pedir o produto, efetuando o pagamento pela internet e retirando o produto no local, seria uma proposta ágil e simples.'
